# Weekly update: insert a new price-report row for
# Femacal de La Calera / Mango just before the existing row 223,
# pushing the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 223 (shifts 223:260 -> 224:261,
# and extends the used range / dimension to A1:T261).
$ws.Rows.Item(223).Insert()

# Fill the newly inserted row 223 with this week's data.
$ws.Cells.Item(223, 1).Value  = 3
$ws.Cells.Item(223, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(223, 3).Value  = "Coquimbo"
$ws.Cells.Item(223, 4).Value  = 44522
$ws.Cells.Item(223, 5).Value  = 5
$ws.Cells.Item(223, 6).Value  = "Fruta"
$ws.Cells.Item(223, 7).Value  = 100108
$ws.Cells.Item(223, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(223, 9).Value  = 100108002
$ws.Cells.Item(223, 10).Value = "Mango"
$ws.Cells.Item(223, 11).Value = "Sin especificar"
$ws.Cells.Item(223, 12).Value = "Primera"
$ws.Cells.Item(223, 13).Value = 125
$ws.Cells.Item(223, 14).Value = 6000
$ws.Cells.Item(223, 15).Value = 6300
$ws.Cells.Item(223, 16).Value = 6156
$ws.Cells.Item(223, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(223, 18).Value = "Perú"
$ws.Cells.Item(223, 19).Value = 1539
$ws.Cells.Item(223, 20).Value = 4
